$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview" (sheet1): drop the d0c8a325 row (row 3), update
# the status text for the 7de87ca3 row, and fix up the hyperlinks.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Rows("3").Delete()

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed903f67c671b4575f21c664cdb5ce855e4f276/e2e/7de87ca3-4966-49b8-8cf0-c1197df20597.md", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ce8a07754b4bff485bfe9c3e155e551d2376b5a/e2e/d0c8a325-d24e-42cb-bdf3-8e1d9193d67f.md", "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "zh-cn" (sheet2): drop the d0c8a325 row (row 3), update the
# status text and the handoff datetime for the 7de87ca3 row, and
# fix up the hyperlinks.
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B2").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "2016-03-10 05:56:38"

$wsZhCn.Rows("3").Delete()

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed903f67c671b4575f21c664cdb5ce855e4f276/e2e/7de87ca3-4966-49b8-8cf0-c1197df20597.md", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0a07df28d28dd490e4db6208ac94739db3c94f20/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.zh-cn.xlf", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/8ee4bceaa3386895ddd26295bffcb9bb5d903e82/e2e/7de87ca3-4966-49b8-8cf0-c1197df20597.md", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/32b2dc9077047e2bb02a28aee6795f53dce8b667/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.zh-cn.xlf", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed903f67c671b4575f21c664cdb5ce855e4f276/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------
# Sheet "de-de" (sheet3): drop the d0c8a325 row (row 3), update the
# status text and the handoff datetime for the 7de87ca3 row, and
# fix up the hyperlinks.
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B2").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "2016-03-10 05:56:47"

$wsDeDe.Rows("3").Delete()

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed903f67c671b4575f21c664cdb5ce855e4f276/e2e/7de87ca3-4966-49b8-8cf0-c1197df20597.md", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8945ba2fc93e0f2c9e9e48f65c9e35cd43e388bc/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.de-de.xlf", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/eb6e9dd65095cbde84a64bd962eadce581455039/e2e/7de87ca3-4966-49b8-8cf0-c1197df20597.md", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1366ad966e7d3bc4290524cfd3d62e49b6c4d255/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.de-de.xlf", "", "", "7de87ca3-4966-49b8-8cf0-c1197df20597.01e140b920f478e5639b9579acb29562e9f4ceb0.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8ed903f67c671b4575f21c664cdb5ce855e4f276/.localization-config", "", "", ".localization-config")

$wb.Save()
